$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"  = 9.264399999999991
    "D4"  = -7.7625
    "D7"  = -7.848700000000004
    "A9"  = -22.0313
    "B9"  = 6.297500000000007
    "C9"  = -11.9732
    "B11" = 5.821700000000002
    "D11" = -7.623100000000002
    "D15" = -8.438799999999997
    "A18" = -22.0407
    "A20" = -19.97269999999998
    "B23" = 8.8901
    "B24" = 6.152200000000001
    "B26" = 4.794400000000001
    "A27" = -22.13980000000002
    "C27" = -13.37019999999999
    "C29" = -11.4732
    "D30" = -7.203499999999996
    "C32" = -12.0162
    "B34" = 9.781700000000008
    "A35" = -22.1462
    "B35" = 4.7599
    "C37" = -12.711
    "C38" = -11.6718
    "D38" = -7.338000000000001
    "D39" = -8.261099999999995
    "C41" = -12.8059
    "D43" = -7.508700000000008
    "C45" = -14.03789999999999
    "D47" = -7.068099999999998
    "B48" = 5.181600000000003
    "B49" = 4.773899999999998
    "C51" = -11.6987
    "B52" = 5.958600000000001
    "C57" = -13.63099999999998
    "C64" = -10.2121
    "B66" = 5.986599999999997
    "B67" = 5.495600000000001
    "A69" = -21.67869999999998
    "D75" = -7.611799999999999
    "A76" = -19.51699999999998
    "A78" = -21.8316
    "B78" = 6.083000000000004
    "B80" = 9.7605
    "A82" = -21.88619999999999
    "C82" = -11.3022
    "A83" = -21.52049999999999
    "D91" = -7.417199999999999
    "D92" = -6.507900000000003
    "A93" = -21.51459999999999
    "C93" = -10.1492
    "D95" = -7.440800000000005
    "D96" = -8.579199999999995
    "B99" = 6.494500000000001
    "C102" = -12.0566
    "B104" = 9.904200000000003
    "C105" = -12.92550000000001
    "D105" = -8.168300000000002
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
